# Grad schools workbook update: add Berkeley and UChicago faculty interests
# and related application-details columns (Fee, Resume?, SoP Notes, Writing
# sample length, # letters of rec, Other requirements?) to the
# "Application Details" sheet, plus a totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Application Details")

# ---------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Deadline "
$ws.Range("D1").Value = "Faculty to mention "
$ws.Range("E1").Value = "Fee"
$ws.Range("F1").Value = "Resume? "
$ws.Range("G1").Value = "SoP Notes"
$ws.Range("H1").Value = "Writing sample length"
$ws.Range("I1").Value = "# letters of rec"
$ws.Range("J1").Value = "Other requirements?"

# ---------------------------------------------------------------------
# Row 2: UChicago / Public Policy PhD
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "12/06/2023"
$ws.Range("D2").Value = "Joshua Gottlieb -- place-based policies. Jeffrey Grogger -- inequality and applied micro. Kelly Hallburg -- segregation, disinvestment. Damon Jones -- racial differences in financial outcomes. Paula Worthington -- metropolitan investment and super applied. "
$ws.Range("E2").Value = 200
$ws.Range("F2").Value = "Y"
$ws.Range("G2").Value = "800 word research statement "
$ws.Range("H2").Value = "TBD"
$ws.Range("I2").Value = 3

# ---------------------------------------------------------------------
# Row 3: UChicago / MSCAPP
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "10/05/2023"
$ws.Range("D3").Formula = "=D2"
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = "Y"
$ws.Range("G3").Value = "300 word motivation statement"
$ws.Range("H3").Value = "N/A"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = "Additional optional essays"

# ---------------------------------------------------------------------
# Row 4: UChicago / MPP
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "10/05/2023"
$ws.Range("D4").Formula = "=D2"
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = "Y"
$ws.Range("G4").Value = "300 word motivation statement"
$ws.Range("H4").Value = "N/A"
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = "Additional optional essays"

# ---------------------------------------------------------------------
# Row 5: Berkeley / Econ PhD
# ---------------------------------------------------------------------
$ws.Range("D5").ClearContents()
$ws.Range("C5").Value = "12/04/2023"
$ws.Range("D5").Value = "Cecile Gaubert -- high-speed rail and spatial sorting & inequality. Hilary Hoynes -- inequality. Pat Kline -- trends in spatial inequality."
$ws.Range("E5").Value = 135
$ws.Range("F5").Value = "Optional "
$ws.Range("G5").Value = "3 pages double spaced"
$ws.Range("H5").Value = "Optional"
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = "Personal history statement (same length as SoP). Must explicitly mention 2 faculty members to work with. "

# ---------------------------------------------------------------------
# Row 6: NYU / Econ PhD
# ---------------------------------------------------------------------
$ws.Range("E6").ClearContents()
$ws.Range("C6").Value = "12/18/2023"
$ws.Range("E6").Value = 110

# ---------------------------------------------------------------------
# Row 9: Georgetown / MS in Data Science for Public Policy
# ---------------------------------------------------------------------
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C9").Value = "12/01/2023"

# ---------------------------------------------------------------------
# Row 25: totals row
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Total applications cost: "
$ws.Range("B25").Formula = "=SUM(E:E)"

Write-Host "Data updated"
